$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 4154.8945
$ws.Range("I19").Value = 4676.2
$ws.Range("J19").Value = 2200
$ws.Range("K19").Value = 4676.2
$ws.Range("L19").Value = 2200
$ws.Range("M19").Value = -4501.2
$ws.Range("N19").Value = -2550

$ws.Range("H33").Value = 152.29033
$ws.Range("I33").Value = 154.17392
$ws.Range("K33").Value = 154.17392
$ws.Range("M33").Value = 74.82607999999999

$ws.Range("H92").Value = 58823948
$ws.Range("I92").Value = 62500384
$ws.Range("K92").Value = 62500384
$ws.Range("M92").Value = -62499136

$ws.Range("H103").Value = 591.1111
$ws.Range("J103").Value = 554.6667
$ws.Range("L103").Value = 1664.0001
$ws.Range("N103").Value = -2836.0001

$ws.Range("H107").Value = 639.85
$ws.Range("I107").Value = 611.2353000000001
$ws.Range("K107").Value = 611.2353000000001
$ws.Range("M107").Value = 1308.7647

$ws.Range("H111").Value = 5343
$ws.Range("I111").Value = 5528.727
$ws.Range("K111").Value = 16586.181
$ws.Range("M111").Value = -13519.181

$ws.Range("H112").Value = 3144.925
$ws.Range("J112").Value = 3144.925
$ws.Range("L112").Value = 9434.775000000001
$ws.Range("N112").Value = -11650.775

$ws.Range("H113").Value = 1500
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 1750
$ws.Range("K113").Value = 1000
$ws.Range("L113").Value = 1750
$ws.Range("M113").Value = 2254
$ws.Range("N113").Value = -8258

$ws.Range("H115").Value = 579.125
$ws.Range("J115").Value = 545.5
$ws.Range("L115").Value = 1636.5
$ws.Range("N115").Value = -4770.5

$ws.Range("H116").Value = 14482429
$ws.Range("I116").Value = 21432076
$ws.Range("J116").Value = 3999.0833
$ws.Range("K116").Value = 21432076
$ws.Range("L116").Value = 3999.0833
$ws.Range("M116").Value = -21428634
$ws.Range("N116").Value = -10883.0833

$ws.Range("H129").Value = 1727.3334
$ws.Range("J129").Value = 1882.6428
$ws.Range("L129").Value = 5647.928400000001
$ws.Range("N129").Value = -15647.9284

$ws.Range("H137").Value = 8590582
$ws.Range("J137").Value = 18527036
$ws.Range("L137").Value = 55581108
$ws.Range("N137").Value = -55586208

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3113.862
$ws.Range("I122").Value = 1714.381
$ws.Range("K122").Value = 5143.143
$ws.Range("M122").Value = -2693.143

$ws.Range("H132").Value = 4205.74
$ws.Range("I132").Value = 1525.6571
$ws.Range("K132").Value = 4576.971299999999
$ws.Range("M132").Value = -2046.971299999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2074.1765
$ws.Range("I94").Value = 2017.5333
$ws.Range("J94").Value = 2499
$ws.Range("K94").Value = 2017.5333
$ws.Range("L94").Value = 2499
$ws.Range("M94").Value = -1566.5333
$ws.Range("N94").Value = -3401

$ws.Range("H99").Value = 4572.75
$ws.Range("I99").Value = 5180.5
$ws.Range("J99").Value = 2749.5
$ws.Range("K99").Value = 5180.5
$ws.Range("L99").Value = 2749.5
$ws.Range("M99").Value = -3682.5
$ws.Range("N99").Value = -5745.5

$ws.Range("H134").Value = 6719.5684
$ws.Range("I134").Value = 3908.15
$ws.Range("J134").Value = 9062.416999999999
$ws.Range("K134").Value = 11724.45
$ws.Range("L134").Value = 27187.251
$ws.Range("M134").Value = -9189.450000000001
$ws.Range("N134").Value = -32257.251

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 35717190
$ws.Range("I31").Value = 41668644
$ws.Range("K31").Value = 41668644
$ws.Range("M31").Value = -41668349

$ws.Range("H34").Value = 35717190
$ws.Range("I34").Value = 41668644
$ws.Range("K34").Value = 41668644
$ws.Range("M34").Value = -41668442

$ws.Range("H99").Value = 7226.263
$ws.Range("I99").Value = 7595.2383
$ws.Range("K99").Value = 7595.2383
$ws.Range("M99").Value = -6097.2383

$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").ClearContents()
$ws.Range("N108").Value = 0

$ws.Range("H126").Value = 7226.263
$ws.Range("I126").Value = 7595.2383
$ws.Range("K126").Value = 22785.7149
$ws.Range("M126").Value = -20315.7149

$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").ClearContents()
$ws.Range("N127").Value = 0

$ws.Range("H132").Value = 121215160
$ws.Range("I132").Value = 148149970
$ws.Range("J132").Value = 8500
$ws.Range("K132").Value = 444449910
$ws.Range("L132").Value = 25500
$ws.Range("M132").Value = -444447380
$ws.Range("N132").Value = -30560

$ws.Range("H134").Value = 1560.381
$ws.Range("I134").Value = 1653.8334
$ws.Range("J134").Value = 999.6667
$ws.Range("K134").Value = 4961.5002
$ws.Range("L134").Value = 2999.0001
$ws.Range("M134").Value = -2426.5002
$ws.Range("N134").Value = -8069.0001

$ws.Range("H137").Value = 106760
$ws.Range("J137").Value = 106760
$ws.Range("L137").Value = 106760
$ws.Range("N137").Value = -116960

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 243.2
$ws.Range("I8").Value = 243.2
$ws.Range("K8").Value = 729.5999999999999
$ws.Range("M8").Value = -590.5999999999999

$ws.Range("H12").Value = 207.92308
$ws.Range("I12").Value = 240.5
$ws.Range("K12").Value = 721.5
$ws.Range("M12").Value = -548.5

$ws.Range("H19").Value = 4000
$ws.Range("J19").Value = 4000
$ws.Range("L19").Value = 12000
$ws.Range("N19").Value = -12348

$ws.Range("H21").Value = 1299.3334
$ws.Range("I21").Value = 1299.3334
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 3898.0002
$ws.Range("L21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -3725.0002

$ws.Range("H25").Value = 834
$ws.Range("I25").Value = 834
$ws.Range("K25").Value = 2502
$ws.Range("M25").Value = -2333

$ws.Range("H29").Value = 63.333332
$ws.Range("I29").Value = 70
$ws.Range("K29").Value = 210
$ws.Range("M29").Value = 67

$ws.Range("H30").Value = 834
$ws.Range("I30").Value = 834
$ws.Range("K30").Value = 2502
$ws.Range("M30").Value = -2400

$ws.Range("H31").Value = 14997.333
$ws.Range("J31").Value = 14997.5
$ws.Range("L31").Value = 44992.5
$ws.Range("N31").Value = -45568.5

$ws.Range("H35").Value = 357.6
$ws.Range("J35").Value = 194.5
$ws.Range("L35").Value = 583.5
$ws.Range("N35").Value = -1159.5

$ws.Range("H36").Value = 1514.7273
$ws.Range("I36").Value = 1160.5
$ws.Range("K36").Value = 3481.5
$ws.Range("M36").Value = -3312.5

$ws.Range("H140").Value = 2064.8333
$ws.Range("I140").Value = 1951
$ws.Range("K140").Value = 5853
$ws.Range("M140").Value = -673

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6405.4
$ws.Range("I70").Value = 6072.1665
$ws.Range("J70").Value = 6627.5557
$ws.Range("K70").Value = 6072.1665
$ws.Range("L70").Value = 6627.5557
$ws.Range("M70").Value = -5802.1665
$ws.Range("N70").Value = -7167.5557

$ws.Range("H73").Value = 6405.4
$ws.Range("I73").Value = 6072.1665
$ws.Range("J73").Value = 6627.5557
$ws.Range("K73").Value = 6072.1665
$ws.Range("L73").Value = 6627.5557
$ws.Range("M73").Value = -5136.1665
$ws.Range("N73").Value = -8499.555700000001

$ws.Range("H80").Value = 67579.17999999999
$ws.Range("I80").Value = 112255.336
$ws.Range("K80").Value = 112255.336
$ws.Range("M80").Value = -111257.336

$ws.Range("H83").Value = 67579.17999999999
$ws.Range("I83").Value = 112255.336
$ws.Range("K83").Value = 561276.6799999999
$ws.Range("M83").Value = -556284.6799999999

$ws.Range("H93").Value = 45000
$ws.Range("J93").Value = 45000
$ws.Range("L93").Value = 45000
$ws.Range("N93").Value = -48744

$ws.Range("H107").Value = 1579.8
$ws.Range("I107").Value = 1224.75
$ws.Range("K107").Value = 1224.75
$ws.Range("M107").Value = 695.25

$ws.Range("H122").Value = 562057
$ws.Range("I122").Value = 3335966
$ws.Range("K122").Value = 10007898
$ws.Range("M122").Value = -10005448

$ws.Range("H126").Value = 2958.361
$ws.Range("I126").Value = 1599.9259
$ws.Range("K126").Value = 4799.7777
$ws.Range("M126").Value = -2329.7777

$ws.Range("H132").Value = 74762.57000000001
$ws.Range("I132").Value = 93909.32000000001
$ws.Range("K132").Value = 281727.96
$ws.Range("M132").Value = -279197.96

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 23814190
$ws.Range("I40").Value = 47621380
$ws.Range("K40").Value = 47621380
$ws.Range("M40").Value = -47621244

$ws.Range("H43").Value = 4898666.5

$ws.Range("H122").Value = 7651.3
$ws.Range("I122").Value = 3740.2
$ws.Range("J122").Value = 11562.4
$ws.Range("K122").Value = 11220.6
$ws.Range("L122").Value = 34687.2
$ws.Range("M122").Value = -8770.599999999999
$ws.Range("N122").Value = -39587.2

$ws.Range("H136").Value = 5142.636
$ws.Range("I136").Value = 3016.3845
$ws.Range("J136").Value = 8213.888999999999
$ws.Range("K136").Value = 9049.1535
$ws.Range("L136").Value = 24641.667
$ws.Range("M136").Value = -6499.1535
$ws.Range("N136").Value = -29741.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7072.85
$ws.Range("I81").Value = 6445.5
$ws.Range("K81").Value = 12891
$ws.Range("M81").Value = -11830

$ws.Range("H84").Value = 7072.85
$ws.Range("I84").Value = 6445.5
$ws.Range("K84").Value = 64455
$ws.Range("M84").Value = -59151
